{"js": "// Load the first paragraph of the body (the one holding the\n// **ID__AFFARS_...__ID** marker) so we can rewrite its formatting/text.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst targetParagraph = paragraphs.items[0];\n\n// The required change touches the paragraph border / indent (pPr) as well\n// as the runs inside it (drop the trailing \" \" run, retarget the ID text).\n// The Word JS API does not expose a way to add a paragraph border whose\n// only attribute is `w:space` (no line style/width/color) nor to delete a\n// single run in one call, so we rebuild the whole paragraph via a flat-OPC\n// OOXML fragment and replace the paragraph's range with it. This keeps\n// every other paragraph/section in the document completely untouched.\nconst paragraphRange = targetParagraph.getRange();\n\nconst flatOpcXml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:pPr>\n              <w:pBdr>\n                <w:top w:space=\"5\"/>\n                <w:left w:space=\"5\"/>\n                <w:bottom w:space=\"5\"/>\n                <w:right w:space=\"5\"/>\n              </w:pBdr>\n              <w:spacing w:after=\"0\"/>\n              <w:ind w:left=\"225\"/>\n              <w:jc w:val=\"left\"/>\n            </w:pPr>\n            <w:r>\n              <w:rPr>\n                <w:rFonts w:ascii=\"Times New Roman\" w:hAnsi=\"Times New Roman\"/>\n                <w:b w:val=\"false\"/>\n                <w:i w:val=\"false\"/>\n                <w:color w:val=\"000000\"/>\n                <w:sz w:val=\"22\"/>\n              </w:rPr>\n              <w:t>**ID__AFFARS_SMC_PGI_5315__ID**</w:t>\n            </w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\nparagraphRange.insertOoxml(flatOpcXml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Target the first paragraph of the document (the hidden \"**ID__...__ID**\"\n# marker paragraph) and apply the requested \"last minute updates\":\n#   1. Replace the ID placeholder text and drop the trailing space run.\n#   2. Bump the paragraph's left indent from 120 to 225 twips (6pt -> 11.25pt).\n#   3. Add a paragraph border whose edges only carry a 5-twip text-to-border\n#      spacing (no line/weight/color), i.e. <w:pBdr><w:top w:space=\"5\"/>...).\n\n$d = $word.ActiveDocument\n$p = $d.Paragraphs(1)\n\n# 1) Collapse the paragraph's two runs (\"**ID__...__ID**\" + trailing \" \")\n#    into a single run holding the new id text, leaving the paragraph mark\n#    (and its formatting) untouched.\n$bodyStart = $p.Range.Start\n$bodyEnd = $p.Range.End - 1\n$bodyRange = $d.Range($bodyStart, $bodyEnd)\n$bodyRange.Text = \"**ID__AFFARS_SMC_PGI_5315__ID**\"\n\n# 2) w:ind w:left=\"120\" -> w:ind w:left=\"225\" (twips / 20 = points)\n$p.Range.ParagraphFormat.LeftIndent = 225 / 20\n\n# 3) w:pBdr with w:space=\"5\" on all four edges\n$p.Range.ParagraphFormat.Borders.DistanceFromTop = 5\n$p.Range.ParagraphFormat.Borders.DistanceFromLeft = 5\n$p.Range.ParagraphFormat.Borders.DistanceFromBottom = 5\n$p.Range.ParagraphFormat.Borders.DistanceFromRight = 5\n"}
